$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")

# New shared strings get interned in the order the values are written, so
# write them in the same order as the target workbook to line up indices.
$ws.Range("B13").Value = "C:/Users/administrator.ENERGIZANDO/Desktop/SISTEMA UNO CLOUD.lnk"
$ws.Range("D12").Value = "Diccionario para la nevagacion de Siesa"
$ws.Range("D13").Value = "Link  a ejecutable de Siesa 8.5"
$ws.Range("A13").Value = "RutaEjecutableSiesa"

# Complete the "Tipo" column for the existing and new rows.
$ws.Range("C12").Value = "Valor"
$ws.Range("C13").Value = "Valor"

# Match author's final cursor/selection position on this sheet.
$ws.Range("B15").Select() | Out-Null
